$d = $word.ActiveDocument

# Locate the placeholder paragraph in the "Alternative Courses" section.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.TrimEnd() -eq "Line n: <alternative course of action>") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $r = $target.Range
    # Exclude the trailing paragraph mark so only the run content is replaced,
    # keeping the paragraph's own properties (numbering, style, ids, etc.) intact.
    $contentRange = $d.Range($r.Start, $r.End - 1)

    $xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512"><pkg:xmlData><Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/></Relationships></pkg:xmlData></pkg:part><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">If there is an issue in generating recommendations, the user will be informed </w:t></w:r><w:r><w:t>about</w:t></w:r><w:r><w:t xml:space="preserve"> the issue.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

    $contentRange.InsertXML($xml)
}
